$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 57114.285
$ws.Range("J133").Value = 57114.285
$ws.Range("L133").Value = 57114.285
$ws.Range("N133").Value = -67234.285

$ws.Range("H137").Value = 1039.0513
$ws.Range("I137").Value = 805.5
$ws.Range("J137").Value = 1506.1538
$ws.Range("K137").Value = 2416.5
$ws.Range("L137").Value = 4518.4614
$ws.Range("M137").Value = 133.5
$ws.Range("N137").Value = -9618.4614

$ws.Range("H141").Value = 2222.4856
$ws.Range("I141").Value = 1816.0667
$ws.Range("J141").Value = 4661
$ws.Range("K141").Value = 5448.2001
$ws.Range("L141").Value = 13983
$ws.Range("M141").Value = -268.2001
$ws.Range("N141").Value = -24343

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9652.219999999999
$ws.Range("I32").Value = 7438.0444
$ws.Range("J32").Value = 29579.8
$ws.Range("K32").Value = 7438.0444
$ws.Range("L32").Value = 29579.8
$ws.Range("M32").Value = -7151.0444
$ws.Range("N32").Value = -30153.8

$ws.Range("H74").Value = 2006.0869
$ws.Range("I74").Value = 2165.0881
$ws.Range("J74").Value = 1555.5834
$ws.Range("K74").Value = 2165.0881
$ws.Range("L74").Value = 1555.5834
$ws.Range("M74").Value = -1291.0881
$ws.Range("N74").Value = -3303.5834

$ws.Range("H77").Value = 2006.0869
$ws.Range("I77").Value = 2165.0881
$ws.Range("J77").Value = 1555.5834
$ws.Range("K77").Value = 10825.4405
$ws.Range("L77").Value = 7777.916999999999
$ws.Range("M77").Value = -6457.440500000001
$ws.Range("N77").Value = -16513.917

$ws.Range("H98").Value = 13538.5
$ws.Range("J98").Value = 13538.5
$ws.Range("L98").Value = 13538.5
$ws.Range("N98").Value = -19528.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1644.3572
$ws.Range("I99").Value = 1233.2963
$ws.Range("J99").Value = 2384.2666
$ws.Range("K99").Value = 1233.2963
$ws.Range("L99").Value = 2384.2666
$ws.Range("M99").Value = 264.7037
$ws.Range("N99").Value = -5380.2666

$ws.Range("H105").Value = 3168.3333
$ws.Range("I105").Value = 2536.6667
$ws.Range("J105").Value = 3800
$ws.Range("K105").Value = 2536.6667
$ws.Range("L105").Value = 3800
$ws.Range("M105").Value = -789.6667000000002
$ws.Range("N105").Value = -7294

$ws.Range("H128").Value = 3991
$ws.Range("I128").Value = 3991
$ws.Range("K128").Value = 11973
$ws.Range("M128").Value = -9483

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2756.1794
$ws.Range("I31").Value = 1492
$ws.Range("J31").Value = 3635.6086
$ws.Range("K31").Value = 1492
$ws.Range("L31").Value = 3635.6086
$ws.Range("M31").Value = -1197
$ws.Range("N31").Value = -4225.6086

$ws.Range("H34").Value = 2756.1794
$ws.Range("I34").Value = 1492
$ws.Range("J34").Value = 3635.6086
$ws.Range("K34").Value = 1492
$ws.Range("L34").Value = 3635.6086
$ws.Range("M34").Value = -1290
$ws.Range("N34").Value = -4039.6086

$ws.Range("H62").Value = 3016.1667
$ws.Range("I62").Value = 2441
$ws.Range("K62").Value = 2441
$ws.Range("M62").Value = -1817

$ws.Range("H65").Value = 3016.1667
$ws.Range("I65").Value = 2441
$ws.Range("K65").Value = 12205
$ws.Range("M65").Value = -9085

$ws.Range("H99").Value = 8187.2856
$ws.Range("I99").Value = 9462.200000000001
$ws.Range("J99").Value = 5000
$ws.Range("K99").Value = 9462.200000000001
$ws.Range("L99").Value = 5000
$ws.Range("M99").Value = -7964.200000000001
$ws.Range("N99").Value = -7996

$ws.Range("H126").Value = 8187.2856
$ws.Range("I126").Value = 9462.200000000001
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 28386.6
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -25916.6
$ws.Range("N126").Value = -19940

$ws.Range("H132").Value = 2005.1666
$ws.Range("I132").Value = 1745.8
$ws.Range("J132").Value = 2134.85
$ws.Range("K132").Value = 5237.4
$ws.Range("L132").Value = 6404.549999999999
$ws.Range("M132").Value = -2707.4
$ws.Range("N132").Value = -11464.55

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1002.64105
$ws.Range("I5").Value = 838.1429000000001
$ws.Range("J5").Value = 1194.5555
$ws.Range("K5").Value = 2514.4287
$ws.Range("L5").Value = 3583.6665
$ws.Range("M5").Value = -2402.4287
$ws.Range("N5").Value = -3807.6665

$ws.Range("H131").Value = 897.1
$ws.Range("I131").Value = 534
$ws.Range("J131").Value = 1092.6154
$ws.Range("K131").Value = 1602
$ws.Range("L131").Value = 3277.8462
$ws.Range("M131").Value = 3438
$ws.Range("N131").Value = -13357.8462

$ws.Range("H135").Value = 1002.64105
$ws.Range("I135").Value = 838.1429000000001
$ws.Range("J135").Value = 1194.5555
$ws.Range("K135").Value = 7543.2861
$ws.Range("L135").Value = 10750.9995
$ws.Range("M135").Value = -5008.2861
$ws.Range("N135").Value = -15820.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 30380.5
$ws.Range("J39").Value = 30380.5
$ws.Range("L39").Value = 30380.5
$ws.Range("N39").Value = -31444.5

$ws.Range("H96").Value = 12000
$ws.Range("J96").Value = 12000
$ws.Range("L96").Value = 12000
$ws.Range("N96").Value = -17492

$ws.Range("H132").Value = 1760.3036
$ws.Range("I132").Value = 1196.7
$ws.Range("J132").Value = 3169.3125
$ws.Range("K132").Value = 3590.1
$ws.Range("L132").Value = 9507.9375
$ws.Range("M132").Value = -1060.1
$ws.Range("N132").Value = -14567.9375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2778494.2
$ws.Range("I22").Value = 4762390
$ws.Range("J22").Value = 1040
$ws.Range("K22").Value = 4762390
$ws.Range("L22").Value = 1040
$ws.Range("M22").Value = -4762095
$ws.Range("N22").Value = -1630

$ws.Range("H27").Value = 2778494.2
$ws.Range("I27").Value = 4762390
$ws.Range("J27").Value = 1040
$ws.Range("K27").Value = 4762390
$ws.Range("L27").Value = 1040
$ws.Range("M27").Value = -4762283
$ws.Range("N27").Value = -1254

$ws.Range("H40").Value = 1655.1765
$ws.Range("I40").Value = 1667
$ws.Range("J40").Value = 1600
$ws.Range("K40").Value = 1667
$ws.Range("L40").Value = 1600
$ws.Range("M40").Value = -1531
$ws.Range("N40").Value = -1872

$ws.Range("H132").Value = 15160938
$ws.Range("I132").Value = 20845018
$ws.Range("J132").Value = 3391.111
$ws.Range("K132").Value = 62535054
$ws.Range("L132").Value = 10173.333
$ws.Range("M132").Value = -62532524
$ws.Range("N132").Value = -15233.333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 30000
$ws.Range("J95").Value = 30000
$ws.Range("L95").Value = 30000
$ws.Range("N95").Value = -35492

$ws.Range("H98").Value = 29000
$ws.Range("J98").Value = 29000
$ws.Range("L98").Value = 29000
$ws.Range("N98").Value = -34990

